$wb = $excel.ActiveWorkbook

# ----- Sheet "Results": update optimization outputs -----
$ws = $wb.Worksheets.Item("Results")

# Row 2 (Bid ID 1, Split A)
$ws.Range("G2").Value = "A"
$ws.Range("H2").Value = 50
$ws.Range("I2").Value = "1%"
$ws.Range("J2").Value = 49.5
$ws.Range("K2").Value = 34650
$ws.Range("M2").Value = 35350
$ws.Range("N2").Value = "10%"
$ws.Range("O2").Value = 3465

# Row 3 (Bid ID 2, Split A)
$ws.Range("G3").Value = "A"
$ws.Range("H3").Value = 70
$ws.Range("I3").Value = "1%"
$ws.Range("J3").Value = 69.3
$ws.Range("K3").Value = 623700
$ws.Range("M3").Value = 780300
$ws.Range("N3").Value = "10%"
$ws.Range("O3").Value = 62370

# Row 4 (Bid ID 3, Split A)
$ws.Range("F4").Value = 253800
$ws.Range("K4").Value = 32670
$ws.Range("L4").Value = 600
$ws.Range("M4").Value = 221130
$ws.Range("N4").Value = "10%"
$ws.Range("O4").Value = 3267

# Row 5 (Bid ID 3 -> 4, Split B -> A)
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "A"
$ws.Range("E5").Value = 453
$ws.Range("F5").Value = 2376891
$ws.Range("G5").Value = "A"
$ws.Range("H5").Value = 23
$ws.Range("I5").Value = "1%"
$ws.Range("J5").Value = 22.77
$ws.Range("K5").Value = 119474.19
$ws.Range("L5").Value = 5247
$ws.Range("M5").Value = 2257416.81
$ws.Range("N5").Value = "10%"
$ws.Range("O5").Value = 11947.419

# Row 6 (Bid ID 4, Split A -> B)
$ws.Range("B6").Value = "B"
$ws.Range("F6").Value = 191619
$ws.Range("G6").Value = "C"
$ws.Range("H6").Value = 24
$ws.Range("I6").Value = "4%"
$ws.Range("J6").Value = 23.04
$ws.Range("K6").Value = 9745.92
$ws.Range("L6").Value = 423
$ws.Range("M6").Value = 181873.08
$ws.Range("N6").Value = "7%"
$ws.Range("O6").Value = 682.2144000000001

# Row 7 (Bid ID 5)
$ws.Range("N7").Value = "7%"
$ws.Range("O7").Value = 72.57600000000001

# Row 8 (Bid ID 6)
$ws.Range("I8").Value = "0%"
$ws.Range("J8").Value = 24
$ws.Range("K8").Value = 5808
$ws.Range("M8").Value = 152218

# Row 9 (Bid ID 7)
$ws.Range("N9").Value = "10%"
$ws.Range("O9").Value = 1511.928

# Row 10 (Bid ID 8)
$ws.Range("I10").Value = "0%"
$ws.Range("J10").Value = 13
$ws.Range("K10").Value = 312
$ws.Range("M10").Value = 10632

# Row 11 (Bid ID 9)
$ws.Range("N11").Value = "7%"
$ws.Range("O11").Value = 202.6752

# Row 12 (Bid ID 10)
$ws.Range("I12").Value = "0%"
$ws.Range("J12").Value = 13
$ws.Range("K12").Value = 169
$ws.Range("M12").Value = 2834

# ----- Sheet "Feasibility Notes": model is now optimal -----
$wsNotes = $wb.Worksheets.Item("Feasibility Notes")
$wsNotes.Range("A2").Value = "Model is optimal."

# ----- Sheet "LP Model": updated constraints (rule fixed) -----
$wsLP = $wb.Worksheets.Item("LP Model")
$lpText = @'
\* Sourcing_with_MultiTier_Rebates_Discounts *\
Minimize
OBJ: S_A + S_B + S_C - rebate_A - rebate_B - rebate_C
Subject To
BaseSpend_A: S0_A - 50 x_A_1 - 64 x_A_10 - 70 x_A_2 - 55 x_A_3 - 23 x_A_4
 - 54 x_A_5 - 42 x_A_6 - 23 x_A_7 - 75 x_A_8 - 97 x_A_9 = 0
BaseSpend_B: S0_B - 60 x_B_1 - 13 x_B_10 - 80 x_B_2 - 65 x_B_3 - 75 x_B_4
 - 34 x_B_5 - 24 x_B_6 - 53 x_B_7 - 13 x_B_8 - 56 x_B_9 = 0
BaseSpend_C: S0_C - 55 x_C_1 - 75 x_C_10 - 75 x_C_2 - 60 x_C_3 - 24 x_C_4
 - 24 x_C_5 - 64 x_C_6 - 86 x_C_7 - 24 x_C_8 - 13 x_C_9 = 0
Capacity_B_Bid_ID_1: x_B_1 <= 100000000
Capacity_B_Bid_ID_10: x_B_10 <= 100000000
Capacity_B_Bid_ID_2: x_B_2 <= 100000000
Capacity_B_Bid_ID_3: x_B_3 <= 100000000
Capacity_B_Bid_ID_4: x_B_4 <= 100000000
Capacity_B_Bid_ID_5: x_B_5 <= 100000000
Capacity_B_Bid_ID_6: x_B_6 <= 100000000
Capacity_B_Bid_ID_7: x_B_7 <= 100000000
Capacity_B_Bid_ID_8: x_B_8 <= 100000000
Capacity_B_Bid_ID_9: x_B_9 <= 100000000
Capacity_C_Bid_ID_1: x_C_1 <= 100000000
Capacity_C_Bid_ID_10: x_C_10 <= 100000000
Capacity_C_Bid_ID_2: x_C_2 <= 100000000
Capacity_C_Bid_ID_3: x_C_3 <= 100000000
Capacity_C_Bid_ID_4: x_C_4 <= 100000000
Capacity_C_Bid_ID_5: x_C_5 <= 100000000
Capacity_C_Bid_ID_6: x_C_6 <= 100000000
Capacity_C_Bid_ID_7: x_C_7 <= 100000000
Capacity_C_Bid_ID_8: x_C_8 <= 100000000
Capacity_C_Bid_ID_9: x_C_9 <= 100000000
Demand_1: x_A_1 + x_B_1 + x_C_1 = 700
Demand_10: x_A_10 + x_B_10 + x_C_10 = 13
Demand_2: x_A_2 + x_B_2 + x_C_2 = 9000
Demand_3: x_A_3 + x_B_3 + x_C_3 = 600
Demand_4: x_A_4 + x_B_4 + x_C_4 = 5670
Demand_5: x_A_5 + x_B_5 + x_C_5 = 45
Demand_6: x_A_6 + x_B_6 + x_C_6 = 242
Demand_7: x_A_7 + x_B_7 + x_C_7 = 664
Demand_8: x_A_8 + x_B_8 + x_C_8 = 24
Demand_9: x_A_9 + x_B_9 + x_C_9 = 232
DiscountTierLower_A_0: d_A - 19400000000 z_discount_A_0 >= -19400000000
DiscountTierLower_A_1: - 0.01 S0_A + d_A - 19400000000 z_discount_A_1
 >= -19400000000
DiscountTierLower_B_0: d_B - 97000000000 z_discount_B_0 >= -97000000000
DiscountTierLower_B_1: - 0.03 S0_B + d_B - 97000000000 z_discount_B_1
 >= -97000000000
DiscountTierLower_C_0: d_C - 97000000000 z_discount_C_0 >= -97000000000
DiscountTierLower_C_1: - 0.04 S0_C + d_C - 97000000000 z_discount_C_1
 >= -97000000000
DiscountTierMax_A_0: x_A_1 + x_A_10 + x_A_2 + x_A_3 + x_A_4 + x_A_5 + x_A_6
 + x_A_7 + x_A_8 + x_A_9 + 19400000000 z_discount_A_0 <= 19400001000
DiscountTierMax_B_0: x_B_1 + x_B_10 + x_B_2 + x_B_3 + x_B_4 + x_B_5 + x_B_6
 + x_B_7 + x_B_8 + x_B_9 + 97000000000 z_discount_B_0 <= 97000000500
DiscountTierMax_C_0: x_C_1 + x_C_10 + x_C_2 + x_C_3 + x_C_4 + x_C_5 + x_C_6
 + x_C_7 + x_C_8 + x_C_9 + 97000000000 z_discount_C_0 <= 97000000500
DiscountTierMin_A_0: x_A_1 + x_A_10 + x_A_2 + x_A_3 + x_A_4 + x_A_5 + x_A_6
 + x_A_7 + x_A_8 + x_A_9 >= 0
DiscountTierMin_A_1: x_A_1 + x_A_10 + x_A_2 + x_A_3 + x_A_4 + x_A_5 + x_A_6
 + x_A_7 + x_A_8 + x_A_9 - 1000 z_discount_A_1 >= 0
DiscountTierMin_B_0: x_B_1 + x_B_10 + x_B_2 + x_B_3 + x_B_4 + x_B_5 + x_B_6
 + x_B_7 + x_B_8 + x_B_9 >= 0
DiscountTierMin_B_1: x_B_1 + x_B_10 + x_B_2 + x_B_3 + x_B_4 + x_B_5 + x_B_6
 + x_B_7 + x_B_8 + x_B_9 - 500 z_discount_B_1 >= 0
DiscountTierMin_C_0: x_C_1 + x_C_10 + x_C_2 + x_C_3 + x_C_4 + x_C_5 + x_C_6
 + x_C_7 + x_C_8 + x_C_9 >= 0
DiscountTierMin_C_1: x_C_1 + x_C_10 + x_C_2 + x_C_3 + x_C_4 + x_C_5 + x_C_6
 + x_C_7 + x_C_8 + x_C_9 - 500 z_discount_C_1 >= 0
DiscountTierSelect_A: z_discount_A_0 + z_discount_A_1 = 1
DiscountTierSelect_B: z_discount_B_0 + z_discount_B_1 = 1
DiscountTierSelect_C: z_discount_C_0 + z_discount_C_1 = 1
DiscountTierUpper_A_0: d_A + 19400000000 z_discount_A_0 <= 19400000000
DiscountTierUpper_A_1: - 0.01 S0_A + d_A + 19400000000 z_discount_A_1
 <= 19400000000
DiscountTierUpper_B_0: d_B + 97000000000 z_discount_B_0 <= 97000000000
DiscountTierUpper_B_1: - 0.03 S0_B + d_B + 97000000000 z_discount_B_1
 <= 97000000000
DiscountTierUpper_C_0: d_C + 97000000000 z_discount_C_0 <= 97000000000
DiscountTierUpper_C_1: - 0.04 S0_C + d_C + 97000000000 z_discount_C_1
 <= 97000000000
EffectiveSpend_A: - S0_A + S_A + d_A = 0
EffectiveSpend_B: - S0_B + S_B + d_B = 0
EffectiveSpend_C: - S0_C + S_C + d_C = 0
RebateTierLower_A_0: rebate_A - 19400000000 y_rebate_A_0 >= -19400000000
RebateTierLower_A_1: - 0.1 S_A + rebate_A - 19400000000 y_rebate_A_1
 >= -19400000000
RebateTierLower_B_0: rebate_B - 97000000000 y_rebate_B_0 >= -97000000000
RebateTierLower_B_1: - 0.05 S_B + rebate_B - 97000000000 y_rebate_B_1
 >= -97000000000
RebateTierLower_C_0: rebate_C - 97000000000 y_rebate_C_0 >= -97000000000
RebateTierLower_C_1: - 0.07 S_C + rebate_C - 97000000000 y_rebate_C_1
 >= -97000000000
RebateTierMax_A_0: x_A_1 + x_A_10 + x_A_2 + x_A_3 + x_A_4 + x_A_5 + x_A_6
 + x_A_7 + x_A_8 + x_A_9 + 19400000000 y_rebate_A_0 <= 19400000500
RebateTierMax_B_0: x_B_1 + x_B_10 + x_B_2 + x_B_3 + x_B_4 + x_B_5 + x_B_6
 + x_B_7 + x_B_8 + x_B_9 + 97000000000 y_rebate_B_0 <= 97000000500
RebateTierMax_C_0: x_C_1 + x_C_10 + x_C_2 + x_C_3 + x_C_4 + x_C_5 + x_C_6
 + x_C_7 + x_C_8 + x_C_9 + 97000000000 y_rebate_C_0 <= 97000000700
RebateTierMin_A_0: x_A_1 + x_A_10 + x_A_2 + x_A_3 + x_A_4 + x_A_5 + x_A_6
 + x_A_7 + x_A_8 + x_A_9 >= 0
RebateTierMin_A_1: x_A_1 + x_A_10 + x_A_2 + x_A_3 + x_A_4 + x_A_5 + x_A_6
 + x_A_7 + x_A_8 + x_A_9 - 500 y_rebate_A_1 >= 0
RebateTierMin_B_0: x_B_1 + x_B_10 + x_B_2 + x_B_3 + x_B_4 + x_B_5 + x_B_6
 + x_B_7 + x_B_8 + x_B_9 >= 0
RebateTierMin_B_1: x_B_1 + x_B_10 + x_B_2 + x_B_3 + x_B_4 + x_B_5 + x_B_6
 + x_B_7 + x_B_8 + x_B_9 - 500 y_rebate_B_1 >= 0
RebateTierMin_C_0: x_C_1 + x_C_10 + x_C_2 + x_C_3 + x_C_4 + x_C_5 + x_C_6
 + x_C_7 + x_C_8 + x_C_9 >= 0
RebateTierMin_C_1: x_C_1 + x_C_10 + x_C_2 + x_C_3 + x_C_4 + x_C_5 + x_C_6
 + x_C_7 + x_C_8 + x_C_9 - 700 y_rebate_C_1 >= 0
RebateTierSelect_A: y_rebate_A_0 + y_rebate_A_1 = 1
RebateTierSelect_B: y_rebate_B_0 + y_rebate_B_1 = 1
RebateTierSelect_C: y_rebate_C_0 + y_rebate_C_1 = 1
RebateTierUpper_A_0: rebate_A + 19400000000 y_rebate_A_0 <= 19400000000
RebateTierUpper_A_1: - 0.1 S_A + rebate_A + 19400000000 y_rebate_A_1
 <= 19400000000
RebateTierUpper_B_0: rebate_B + 97000000000 y_rebate_B_0 <= 97000000000
RebateTierUpper_B_1: - 0.05 S_B + rebate_B + 97000000000 y_rebate_B_1
 <= 97000000000
RebateTierUpper_C_0: rebate_C + 97000000000 y_rebate_C_0 <= 97000000000
RebateTierUpper_C_1: - 0.07 S_C + rebate_C + 97000000000 y_rebate_C_1
 <= 97000000000
Transition_10_A: - 13 T_10_A + x_A_10 <= 0
Transition_10_B: - 13 T_10_B + x_B_10 <= 0
Transition_1_B: - 700 T_1_B + x_B_1 <= 0
Transition_1_C: - 700 T_1_C + x_C_1 <= 0
Transition_2_A: - 9000 T_2_A + x_A_2 <= 0
Transition_2_C: - 9000 T_2_C + x_C_2 <= 0
Transition_3_A: - 600 T_3_A + x_A_3 <= 0
Transition_3_B: - 600 T_3_B + x_B_3 <= 0
Transition_4_A: - 5670 T_4_A + x_A_4 <= 0
Transition_4_B: - 5670 T_4_B + x_B_4 <= 0
Transition_5_A: - 45 T_5_A + x_A_5 <= 0
Transition_5_B: - 45 T_5_B + x_B_5 <= 0
Transition_6_A: - 242 T_6_A + x_A_6 <= 0
Transition_6_B: - 242 T_6_B + x_B_6 <= 0
Transition_7_A: - 664 T_7_A + x_A_7 <= 0
Transition_7_B: - 664 T_7_B + x_B_7 <= 0
Transition_8_A: - 24 T_8_A + x_A_8 <= 0
Transition_8_B: - 24 T_8_B + x_B_8 <= 0
Transition_9_A: - 232 T_9_A + x_A_9 <= 0
Transition_9_B: - 232 T_9_B + x_B_9 <= 0
VolAwarded_Agg_0_incumbent_UB: x_A_1 + x_C_10 + x_C_3 + x_C_4 + x_C_5 + x_C_9
 <= 3000
Volume_A: V_A - x_A_1 - x_A_10 - x_A_2 - x_A_3 - x_A_4 - x_A_5 - x_A_6 - x_A_7
 - x_A_8 - x_A_9 = 0
Volume_B: V_B - x_B_1 - x_B_10 - x_B_2 - x_B_3 - x_B_4 - x_B_5 - x_B_6 - x_B_7
 - x_B_8 - x_B_9 = 0
Volume_C: V_C - x_C_1 - x_C_10 - x_C_2 - x_C_3 - x_C_4 - x_C_5 - x_C_6 - x_C_7
 - x_C_8 - x_C_9 = 0
Binaries
T_10_A
T_10_B
T_1_B
T_1_C
T_2_A
T_2_C
T_3_A
T_3_B
T_4_A
T_4_B
T_5_A
T_5_B
T_6_A
T_6_B
T_7_A
T_7_B
T_8_A
T_8_B
T_9_A
T_9_B
y_rebate_A_0
y_rebate_A_1
y_rebate_B_0
y_rebate_B_1
y_rebate_C_0
y_rebate_C_1
z_discount_A_0
z_discount_A_1
z_discount_B_0
z_discount_B_1
z_discount_C_0
z_discount_C_1
End

'@
$wsLP.Range("A2").Value = $lpText
